$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10 (pushes "fossil_routes" and everything below it
# down by one row), to make room for the new "chemical_recycling_pyrolysis"
# parameter right after "chemical_recycling_gasification".
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
